$d = $word.ActiveDocument

# The "_GoBack" bookmark currently wraps the end of the "Float property in
# css(img)" paragraph. In the target document it instead wraps the end of
# the new final paragraph ("Script tag"), so drop it here and we will
# re-add it on the new last paragraph below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Insert a fresh paragraph right after paragraph 6 ("Float property in
# css(img)") to use as an anchor, then replace its content with all of the
# new paragraphs (as literal WordprocessingML) in one shot so that run
# boundaries, the paragraph border, and the bookmark come out exactly as
# authored.
$anchor = $d.Paragraphs(6).Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$anchor.Collapse(0)

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Upwork</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Clear both</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Animate </w:t></w:r><w:r><w:t>.css</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Prompt()</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Google fonts</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Boxmodel</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>T</w:t></w:r><w:r><w:t>ranisition</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Flexbox</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Mailchimp</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Prompt()</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Number(first) for con</w:t></w:r><w:r><w:t>verting it into number</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>A</w:t></w:r><w:r><w:t>lert</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:bottom w:val="double" w:sz="6" w:space="1" w:color="auto"/></w:pBdr></w:pPr><w:r><w:t xml:space="preserve">Creative </w:t></w:r><w:r><w:t>tim</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Script tag</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$target = $d.Paragraphs(7).Range
$target.InsertXML($newParagraphsXml)

# Remove the two trailing paragraphs left over from the original document
# (a lone space, then an empty paragraph) now that the new content sits
# right after "Float property in css(img)" and before them.
$count = $d.Paragraphs.Count
$trailingStart = $d.Paragraphs($count - 1).Range.Start
$trailingEnd = $d.Paragraphs($count).Range.End
$d.Range($trailingStart, $trailingEnd).Delete()
